{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the credit-line paragraph (\"Rasmus Tilljander - rati10@student.bth.se\")\n// and the work-summary paragraph (ends with \"...screen handler.\") by their text\n// content, rather than relying on fixed paragraph indices.\nlet creditPara = null;\nlet summaryPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Rasmus\") !== -1 && t.indexOf(\"Tilljander\") !== -1) {\n    creditPara = paragraphs.items[i];\n  }\n  if (t.indexOf(\"screen handler\") !== -1) {\n    summaryPara = paragraphs.items[i];\n  }\n}\n\nif (!creditPara || !summaryPara) {\n  throw new Error(\"Could not locate target paragraphs\");\n}\n\n// 1) Rebuild the credit-line paragraph so \"Rasmus\" and \"Tilljander\" each sit in\n// their own spell-check-marked run, \"Rasmus\" and the trailing space become\n// separate runs, and every run (plus the paragraph mark) is tagged en-US -\n// matching what Word produces when the paragraph's proofing language is set.\nconst creditOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:pPr><w:pStyle w:val=\"Footer\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Rasmus</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Tilljander</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> - rati10@student.bth.se</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\ncreditPara.insertOoxml(creditOoxml, \"Replace\");\n\n// 2) Append the new sentence about starting the project skeleton to the work\n// summary paragraph, tagged en-US like the rest of the paragraph.\nconst newRange = summaryPara.insertText(\n  \" We also started programming the skeleton for the project.\",\n  \"End\"\n);\nnewRange.languageId = \"en-US\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the credit-line paragraph (\"Rasmus Tilljander - rati10@student.bth.se\")\n# and the work-summary paragraph (ends with \"...screen handler.\") by their text,\n# instead of relying on fixed paragraph indices.\n$creditPara = $null\n$summaryPara = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Rasmus*\" -and $t -like \"*Tilljander*\") {\n        $creditPara = $p\n    }\n    if ($t -like \"*screen handler*\") {\n        $summaryPara = $p\n    }\n}\n\n# 1) Append the new sentence to the end of the work-summary paragraph's text\n# (before its paragraph mark) and tag it en-US like the rest of the paragraph.\n$summaryRange = $summaryPara.Range\n$summaryRange.End = $summaryRange.End - 1\n$summaryRange.Collapse(0)\n$summaryRange.InsertAfter(\" We also started programming the skeleton for the project.\")\n$summaryRange.LanguageID = \"en-US\"\n\n# 2) Rebuild the credit-line paragraph so \"Rasmus\" and \"Tilljander\" each sit in\n# their own spell-check-marked run, \"Rasmus\" and the following space are split\n# into separate runs, and every run (plus the paragraph mark) is tagged en-US -\n# matching what Word produces when the paragraph's proofing language is set.\n$creditRange = $creditPara.Range\n$creditXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p>' +\n    '<w:pPr><w:pStyle w:val=\"Footer\"/><w:rPr><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Rasmus</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t>Tilljander</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> - rati10@student.bth.se</w:t></w:r>' +\n    '</w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n$creditRange.InsertXML($creditXml) | Out-Null\n"}
